# Use exogenous Carbon Intensity Ratios so LCFS works as a biofuel share
# requirement (PEI != 0 for biodiesel in BAU)

$wb = $excel.ActiveWorkbook

$wsCIRbTF = $wb.Worksheets.Item("CIRbTF")

# Data change: biofuel diesel ratio/flag goes from -1 (calculate
# automatically) to 0 (exogenous / overridden value) so that LCFS acts as
# a biofuel share requirement.
$wsCIRbTF.Range("B7").Value = 0

# Make CIRbTF the active / selected sheet (matches tabSelected moving from
# About to CIRbTF, and the new active selection on CIRbTF).
$wsCIRbTF.Activate()
$wsCIRbTF.Range("B8").Select()

$wb.Save()
